# Scheduled-runner market-data refresh for the Leve profitability sheets.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) with
# freshly pulled Universalis price data -- item/recipe columns (A:G) are untouched.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws_ALC = $wb.Worksheets.Item("ALC")
# Row 4: Root Rush / Growth Formula Alpha
$ws_ALC.Range("H4").Value = 333.33334
$ws_ALC.Range("I4").Value = 333.33334
$ws_ALC.Range("K4").Value = 333.33334
$ws_ALC.Range("M4").Value = -219.33334
# Row 6: Days of Chunder / Antidote
$ws_ALC.Range("H6").Value = 1771.5
$ws_ALC.Range("I6").Value = 1771.5
$ws_ALC.Range("K6").Value = 5314.5
$ws_ALC.Range("M6").Value = -5202.5
# Row 11: Gotta Bounce / Rubber
$ws_ALC.Range("H11").Value = 143.44444
$ws_ALC.Range("I11").Value = 143.44444
$ws_ALC.Range("K11").Value = 143.44444
$ws_ALC.Range("M11").Value = -3.444439999999986
# Row 15: Morning Glass of Ether / Ether
$ws_ALC.Range("H15").Value = 1811.4231
$ws_ALC.Range("I15").Value = 1811.4231
$ws_ALC.Range("K15").Value = 5434.2693
$ws_ALC.Range("M15").Value = -5265.2693
# Row 33: Glazed and Confused / Clear Glass Lens
$ws_ALC.Range("H33").Value = 398.84616
$ws_ALC.Range("I33").Value = 390.41666
$ws_ALC.Range("J33").Value = 500
$ws_ALC.Range("K33").Value = 390.41666
$ws_ALC.Range("L33").Value = 500
$ws_ALC.Range("M33").Value = -161.41666
$ws_ALC.Range("N33").Value = -958
# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws_ALC.Range("H86").Value = 4973.5
$ws_ALC.Range("I86").Value = 4973
$ws_ALC.Range("J86").Value = 4974
$ws_ALC.Range("K86").Value = 4973
$ws_ALC.Range("L86").Value = 4974
$ws_ALC.Range("M86").Value = -3850
$ws_ALC.Range("N86").Value = -7220
# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws_ALC.Range("H89").Value = 4973.5
$ws_ALC.Range("I89").Value = 4973
$ws_ALC.Range("J89").Value = 4974
$ws_ALC.Range("K89").Value = 24865
$ws_ALC.Range("L89").Value = 24870
$ws_ALC.Range("M89").Value = -19249
$ws_ALC.Range("N89").Value = -36102
# Row 92: Whinier than the Sword / Enchanted Koppranickel Ink
$ws_ALC.Range("H92").Value = 758.9167
$ws_ALC.Range("I92").Value = 559
$ws_ALC.Range("J92").Value = 1358.6666
$ws_ALC.Range("K92").Value = 559
$ws_ALC.Range("L92").Value = 1358.6666
$ws_ALC.Range("M92").Value = 689
$ws_ALC.Range("N92").Value = -3854.6666
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws_ALC.Range("H98").Value = 793.44446
$ws_ALC.Range("I98").Value = 767.625
$ws_ALC.Range("J98").Value = 1000
$ws_ALC.Range("K98").Value = 767.625
$ws_ALC.Range("L98").Value = 1000
$ws_ALC.Range("M98").Value = 730.375
$ws_ALC.Range("N98").Value = -3996
# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws_ALC.Range("H107").Value = 458.9091
$ws_ALC.Range("I107").Value = 607.1429000000001
$ws_ALC.Range("J107").Value = 199.5
$ws_ALC.Range("K107").Value = 607.1429000000001
$ws_ALC.Range("L107").Value = 199.5
$ws_ALC.Range("M107").Value = 1312.8571
$ws_ALC.Range("N107").Value = -4039.5
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws_ALC.Range("H122").Value = 793.44446
$ws_ALC.Range("I122").Value = 767.625
$ws_ALC.Range("J122").Value = 1000
$ws_ALC.Range("K122").Value = 2302.875
$ws_ALC.Range("L122").Value = 3000
$ws_ALC.Range("M122").Value = 147.125
$ws_ALC.Range("N122").Value = -7900
# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws_ALC.Range("H135").Value = 1378.5714
$ws_ALC.Range("I135").Value = 950
$ws_ALC.Range("K135").Value = 8550
$ws_ALC.Range("M135").Value = -6015
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws_ALC.Range("H138").Value = 4336.931
$ws_ALC.Range("J138").Value = 4527
$ws_ALC.Range("L138").Value = 13581
$ws_ALC.Range("N138").Value = -23861

# --- ARM sheet ---
$ws_ARM = $wb.Worksheets.Item("ARM")
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws_ARM.Range("H132").Value = 2843.9565
$ws_ARM.Range("I132").Value = 2094.3333
$ws_ARM.Range("K132").Value = 6282.999899999999
$ws_ARM.Range("M132").Value = -3752.999899999999

# --- BSM sheet ---
$ws_BSM = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight / Mythrite Nugget
$ws_BSM.Range("H64").Value = 1185.8
$ws_BSM.Range("J64").Value = 1232.25
$ws_BSM.Range("L64").Value = 1232.25
$ws_BSM.Range("N64").Value = -1682.25
# Row 67: Bearing the Brunt (L) / Mythrite Nugget
$ws_BSM.Range("H67").Value = 1185.8
$ws_BSM.Range("J67").Value = 1232.25
$ws_BSM.Range("L67").Value = 1232.25
$ws_BSM.Range("N67").Value = -2792.25
# Row 107: The Gold Experience / Deepgold Nugget
$ws_BSM.Range("H107").Value = 963.1429000000001
$ws_BSM.Range("I107").Value = 871.9091
$ws_BSM.Range("K107").Value = 871.9091
$ws_BSM.Range("M107").Value = 1048.0909
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws_BSM.Range("H134").Value = 2431.5151
$ws_BSM.Range("I134").Value = 2041
$ws_BSM.Range("K134").Value = 6123
$ws_BSM.Range("M134").Value = -3588

# --- CRP sheet ---
$ws_CRP = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws_CRP.Range("H16").Value = 1050
$ws_CRP.Range("I16").Value = 970.2
$ws_CRP.Range("J16").Value = 1149.75
$ws_CRP.Range("K16").Value = 970.2
$ws_CRP.Range("L16").Value = 1149.75
$ws_CRP.Range("M16").Value = -683.2
$ws_CRP.Range("N16").Value = -1723.75
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws_CRP.Range("H58").Value = 2932.5557
$ws_CRP.Range("I58").Value = 2862.4443
$ws_CRP.Range("J58").Value = 3002.6667
$ws_CRP.Range("K58").Value = 2862.4443
$ws_CRP.Range("L58").Value = 3002.6667
$ws_CRP.Range("M58").Value = -2659.4443
$ws_CRP.Range("N58").Value = -3408.6667
# Row 86: Birch, Please / Birch Lumber
$ws_CRP.Range("H86").Value = 21873.6
$ws_CRP.Range("I86").Value = 10033.818
$ws_CRP.Range("K86").Value = 10033.818
$ws_CRP.Range("M86").Value = -8910.817999999999
# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws_CRP.Range("H89").Value = 21873.6
$ws_CRP.Range("I89").Value = 10033.818
$ws_CRP.Range("K89").Value = 50169.09
$ws_CRP.Range("M89").Value = -44553.09
# Row 113: Patient Patients / White Ash Lumber
$ws_CRP.Range("H113").Value = 1050
$ws_CRP.Range("I113").Value = 970.2
$ws_CRP.Range("J113").Value = 1149.75
$ws_CRP.Range("K113").Value = 970.2
$ws_CRP.Range("L113").Value = 1149.75
$ws_CRP.Range("M113").Value = 1199.8
$ws_CRP.Range("N113").Value = -5489.75
# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws_CRP.Range("H122").Value = 1728
$ws_CRP.Range("I122").Value = 2012
$ws_CRP.Range("J122").Value = 1444
$ws_CRP.Range("K122").Value = 6036
$ws_CRP.Range("L122").Value = 4332
$ws_CRP.Range("M122").Value = -3586
$ws_CRP.Range("N122").Value = -9232
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws_CRP.Range("H132").Value = 4445.5625
$ws_CRP.Range("I132").Value = 4331.7
$ws_CRP.Range("K132").Value = 12995.1
$ws_CRP.Range("M132").Value = -10465.1
# Row 136: Turali Quality / Dark Mahogany Lumber
$ws_CRP.Range("H136").Value = 2932.5557
$ws_CRP.Range("I136").Value = 2862.4443
$ws_CRP.Range("J136").Value = 3002.6667
$ws_CRP.Range("K136").Value = 8587.332900000001
$ws_CRP.Range("L136").Value = 9008.000100000001
$ws_CRP.Range("M136").Value = -6037.332900000001
$ws_CRP.Range("N136").Value = -14108.0001
# Row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws_CRP.Range("H141").Value = 51191.23
$ws_CRP.Range("J141").Value = 51191.23
$ws_CRP.Range("L141").Value = 51191.23
$ws_CRP.Range("N141").Value = -61551.23

# --- CUL sheet ---
$ws_CUL = $wb.Worksheets.Item("CUL")
# Row 6: Meat-lover's Special / Marmot Steak
$ws_CUL.Range("H6").Value = 899.75
$ws_CUL.Range("I6").Value = 199.66667
$ws_CUL.Range("J6").Value = 3000
$ws_CUL.Range("K6").Value = 599.00001
$ws_CUL.Range("L6").Value = 9000
$ws_CUL.Range("M6").Value = -486.00001
$ws_CUL.Range("N6").Value = -9226
# Row 13: Fishy Revelations / Braised Pipira
$ws_CUL.Range("H13").Value = 699.5
$ws_CUL.Range("I13").Value = 699.5
$ws_CUL.Range("K13").Value = 2098.5
$ws_CUL.Range("M13").Value = -1930.5
# Row 17: Chew the Fat / Grilled Dodo
$ws_CUL.Range("H17").Value = 1442.5
$ws_CUL.Range("I17").Value = 90
$ws_CUL.Range("J17").Value = 1893.3334
$ws_CUL.Range("K17").Value = 270
$ws_CUL.Range("L17").Value = 5680.0002
$ws_CUL.Range("M17").Value = -101
$ws_CUL.Range("N17").Value = -6018.0002
# Row 107: Slippery Service / Frantoio Oil
$ws_CUL.Range("H107").Value = 1503.625
$ws_CUL.Range("J107").Value = 289.85715
$ws_CUL.Range("L107").Value = 869.5714499999999
$ws_CUL.Range("N107").Value = -4709.571449999999
# Row 117: A Good Omen / Peppered Popotoes
$ws_CUL.Range("H117").Value = 2300
$ws_CUL.Range("J117").Value = 2300
$ws_CUL.Range("L117").Value = 6900
$ws_CUL.Range("N117").Value = -13784
# Row 129: Comfort Food / Yakow Moussaka
$ws_CUL.Range("H129").Value = 1447.75
$ws_CUL.Range("I129").Value = 1286.6
$ws_CUL.Range("J129").Value = 1562.8572
$ws_CUL.Range("K129").Value = 3859.8
$ws_CUL.Range("L129").Value = 4688.571599999999
$ws_CUL.Range("M129").Value = 1140.2
$ws_CUL.Range("N129").Value = -14688.5716

# --- GSM sheet ---
$ws_GSM = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws_GSM.Range("H2").Value = 49
$ws_GSM.Range("J2").Value = 25.75
$ws_GSM.Range("L2").Value = 25.75
$ws_GSM.Range("N2").Value = -251.75
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws_GSM.Range("H126").Value = 1695.5
$ws_GSM.Range("I126").Value = 1256.5
$ws_GSM.Range("J126").Value = 2134.5
$ws_GSM.Range("K126").Value = 3769.5
$ws_GSM.Range("L126").Value = 6403.5
$ws_GSM.Range("M126").Value = -1299.5
$ws_GSM.Range("N126").Value = -11343.5

# --- LTW sheet ---
$ws_LTW = $wb.Worksheets.Item("LTW")
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws_LTW.Range("H68").Value = 2499.25
$ws_LTW.Range("I68").Value = 2499.25
$ws_LTW.Range("K68").Value = 2499.25
$ws_LTW.Range("M68").Value = -1750.25
# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws_LTW.Range("H71").Value = 2499.25
$ws_LTW.Range("I71").Value = 2499.25
$ws_LTW.Range("K71").Value = 12496.25
$ws_LTW.Range("M71").Value = -8752.25

# --- WVR sheet ---
$ws_WVR = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display / Ruby Cotton Cloth
$ws_WVR.Range("H96").Value = 1296.25
$ws_WVR.Range("I96").Value = 1296.25
$ws_WVR.Range("K96").Value = 1296.25
$ws_WVR.Range("M96").Value = 76.75
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws_WVR.Range("H136").Value = 3537.08
$ws_WVR.Range("I136").Value = 3371.3125
$ws_WVR.Range("K136").Value = 10113.9375
$ws_WVR.Range("M136").Value = -7563.9375
